$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the inner variable name values in literal double quotes
$ws.Range("D2").Value = '"PMax"'
$ws.Range("D3").Value = '"KPO4"'

# Update the active selection to match the saved view state
$ws.Range("D4").Select()
